# #5 Notes should be imported. Ensure zipcode format.
# Improvement for mailing address displaying format.
#
# Splits the combined "Address Line 2" (City, State Zip) into three
# discrete columns - City, State, Zipcode - inserted between
# "Address Line 2" and "Phone". Also restores a missed Notes entry and
# tidies the sheet's print setup / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (City, State, Zipcode) right before the
# existing "Phone" column (F). Everything to the right (Phone, Phone Tag,
# Employment Status, Notes, Relationships) shifts from F:J to I:M.
$ws.Columns("F:H").Insert()

# New column headers. (Columns("F:H").Insert() already carried the bold
# header style over from column E, matching the rest of row 1.)
$ws.Range("F1").Value = "City"
$ws.Range("G1").Value = "State"
$ws.Range("H1").Value = "Zipcode"

# Give the new columns a sensible display width, matching the other
# address column.
$ws.Range("F1:H1").ColumnWidth = 16.998697916666668

# Per-row City / State / Zipcode values, derived from the pre-existing
# "Austin, TX 78701" / "New York, NY 80001" style Address Line 2 values.
$cities  = @("Austin ", "Austin ", "Austin ", "Austin ", "Austin ", "Austin ", "Austin ", "New York", "New York", "New York", "New York")
$states  = @("TX", "TX", "TX", "TX", "TX", "TX", "TX", "NY", "NY", "NY", "NY")
$zips    = @(78701, 78702, 78703, 78704, 78705, 78706, 78707, 80001, 80002, 80003, 80004)

for ($i = 0; $i -lt $cities.Length; $i++) {
    $row = 2 + $i
    $ws.Range("F$row").Value = $cities[$i]
    $ws.Range("G$row").Value = $states[$i]
    $ws.Range("H$row").Value = $zips[$i]
}

# Restore the missing note for Bill Clinton (row 5), now in the
# "Notes" column, which moved from I to L.
$ws.Range("L5").Value = "Mr President"

# Print setup tweak for the new, wider mailing-address layout.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Restore the (arbitrary) last-used selection.
$null = $ws.Range("K26").Select()
